$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Refresh the "datetimeFigureOut" date placeholders (slide master,
#    every slide layout, and the notes master) from 1/7/2017 -> 3/15/2017.
# ---------------------------------------------------------------------
function Update-DatePlaceholders($container) {
    for ($j = 1; $j -le $container.Shapes.Count; $j++) {
        $sh = $container.Shapes.Item($j)
        if ($sh.HasTextFrame) {
            if ($sh.TextFrame.TextRange.Text -eq "1/7/2017") {
                $sh.TextFrame.TextRange.Text = "3/15/2017"
            }
        }
    }
}

Update-DatePlaceholders($p.SlideMaster)
Update-DatePlaceholders($p.NotesMaster)

$layouts = $p.SlideMaster.CustomLayouts
for ($i = 1; $i -le $layouts.Count; $i++) {
    Update-DatePlaceholders($layouts.Item($i))
}

# ---------------------------------------------------------------------
# 2) Rename the AddressBook-era labels on slide 1 to the Task-based
#    names ("PersonListPanel" -> "TaskListPanel", "PersonCard" -> "TaskCard").
# ---------------------------------------------------------------------
$s = $p.Slides.Item(1)
for ($k = 1; $k -le $s.Shapes.Count; $k++) {
    $shape = $s.Shapes.Item($k)
    if ($shape.HasTextFrame) {
        $t = $shape.TextFrame.TextRange.Text
        if ($t -eq "PersonListPanel") {
            $shape.TextFrame.TextRange.Text = "TaskListPanel"
        } elseif ($t -eq "PersonCard") {
            $shape.TextFrame.TextRange.Text = "TaskCard"
        }
    }
}
